$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 260
$ws1.Range("F3").Value = 560
$ws1.Range("G5").Value = 49.9
$ws1.Range("F6").Value = 1064
$ws1.Range("F7").Value = 1398
$ws1.Range("F8").Value = 577
$ws1.Range("F9").Value = 97
$ws1.Range("F10").Value = 733
$ws1.Range("F13").Value = 114
$ws1.Range("F15").Value = 1277
$ws1.Range("F16").Value = 95
$ws1.Range("F17").Value = 77
$ws1.Range("F20").Value = 631
$ws1.Range("F21").Value = 27
$ws1.Range("F22").Value = 188
$ws1.Range("F23").Value = 5515
$ws1.Range("F28").Value = 13956
$ws1.Range("F30").Value = 185
$ws1.Range("F31").Value = 85
$ws1.Range("F33").Value = 402
$ws1.Range("F34").Value = 550
$ws1.Range("F35").Value = 4147
$ws1.Range("F36").Value = 89
$ws1.Range("F37").Value = 348

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 260
$ws4.Range("F3").Value = 560
$ws4.Range("G5").Value = 49.9
$ws4.Range("F6").Value = 1064
$ws4.Range("F7").Value = 1398
$ws4.Range("F8").Value = 577
$ws4.Range("F9").Value = 97
$ws4.Range("F10").Value = 733
$ws4.Range("F13").Value = 114
$ws4.Range("F15").Value = 1277
$ws4.Range("F16").Value = 95
$ws4.Range("F17").Value = 77
$ws4.Range("F21").Value = 631
$ws4.Range("F23").Value = 27
$ws4.Range("F24").Value = 188
$ws4.Range("F26").Value = 5515
$ws4.Range("F31").Value = 13956
$ws4.Range("F33").Value = 185
$ws4.Range("F34").Value = 85
$ws4.Range("F36").Value = 402
$ws4.Range("F37").Value = 550
$ws4.Range("F38").Value = 4147
$ws4.Range("F39").Value = 89
$ws4.Range("F40").Value = 348
